$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure columns B:E are treated as Text so numeric-looking strings
# (e.g. "1.00", "94.279.53") are preserved as text, matching the
# original inline-string cell types instead of being auto-converted
# to numbers by Excel's smart entry.
$dataRange = $ws.Range("B2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = '94.279.53'
$ws.Range("E2").Value = '  -1.78%  '
$ws.Range("D3").Value = '3.330.06'
$ws.Range("E3").Value = '  -4.37%  '
$ws.Range("E4").Value = '  +0.12%  '
$ws.Range("D5").Value = '231.74'
$ws.Range("E5").Value = '  -4.62%  '
$ws.Range("D6").Value = '618.85'
$ws.Range("E6").Value = '  -4.47%  '
$ws.Range("D7").Value = '1.38'
$ws.Range("E7").Value = '  -6.07%  '
$ws.Range("D8").Value = '0.388'
$ws.Range("E8").Value = '  -6.55%  '
$ws.Range("D9").Value = '1.00'
$ws.Range("E9").Value = '  +0.02%  '
$ws.Range("D10").Value = '0.935'
$ws.Range("E10").Value = '  -6.34%  '
$ws.Range("D11").Value = '3.327.47'
$ws.Range("E11").Value = '  -4.56%  '
$ws.Range("D12").Value = '42.01'
$ws.Range("E12").Value = '  -2.03%  '
$ws.Range("E13").Value = '  -3.16%  '
$ws.Range("D14").Value = '94.138.70'
$ws.Range("E14").Value = '  -1.63%  '
$ws.Range("D15").Value = '5.95'
$ws.Range("E15").Value = '  -3.31%  '
$ws.Range("D16").Value = '3.959.06'
$ws.Range("E16").Value = '  -4.05%  '
$ws.Range("D17").Value = '0.0000244'
$ws.Range("E17").Value = '  -4.57%  '
$ws.Range("D18").Value = '8.10'
$ws.Range("E18").Value = '  -4.66%  '
$ws.Range("D19").Value = '3.328.18'
$ws.Range("E19").Value = '  -4.27%  '
$ws.Range("D20").Value = '17.25'
$ws.Range("E20").Value = '  -6.05%  '
$ws.Range("D21").Value = '10.96'
$ws.Range("E21").Value = '  -7.74%  '
$ws.Range("D22").Value = '3.49'
$ws.Range("E22").Value = '  +7.12%  '
$ws.Range("D23").Value = '493.90'
$ws.Range("E23").Value = '  -3.58%  '
$ws.Range("D24").Value = '0.453'
$ws.Range("E24").Value = '  -13.01%  '
$ws.Range("D25").Value = '0.0000182'
$ws.Range("E25").Value = '  -6.12%  '
$ws.Range("D26").Value = '6.07'
$ws.Range("E26").Value = '  -9.16%  '
$ws.Range("D27").Value = '89.93'
$ws.Range("E27").Value = '  -2.70%  '
$ws.Range("D28").Value = '11.72'
$ws.Range("E28").Value = '  -5.14%  '
$ws.Range("D29").Value = '3.508.45'
$ws.Range("E29").Value = '  -4.12%  '
$ws.Range("E30").Value = '  +0.07%  '
$ws.Range("D31").Value = '11.12'
$ws.Range("E31").Value = '  -6.64%  '
$ws.Range("D32").Value = '0.137'
$ws.Range("E32").Value = '  -0.78%  '
$ws.Range("D33").Value = '2.63'
$ws.Range("E33").Value = '  -5.21%  '
$ws.Range("D34").Value = '1.00'
$ws.Range("E34").Value = '  +0.50%  '
$ws.Range("D35").Value = '0.174'
$ws.Range("E35").Value = '  -6.02%  '
$ws.Range("D36").Value = '28.36'
$ws.Range("E36").Value = '  -8.84%  '
$ws.Range("D37").Value = '0.531'
$ws.Range("E37").Value = '  -8.19%  '
$ws.Range("D38").Value = '530.22'
$ws.Range("E38").Value = '  +1.57%  '
$ws.Range("B39").Value = 'USDe'
$ws.Range("C39").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D39").Value = '1.00'
$ws.Range("E39").Value = '  +0.00%  '
$ws.Range("B40").Value = 'RenderToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D40").Value = '7.36'
$ws.Range("E40").Value = '  -6.13%  '
$ws.Range("B41").Value = 'Kaspa'
$ws.Range("C41").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D41").Value = '0.148'
$ws.Range("E41").Value = '  -2.50%  '
$ws.Range("B42").Value = 'Fetch.AI'
$ws.Range("C42").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D42").Value = '1.36'
$ws.Range("E42").Value = '  -6.71%  '
$ws.Range("D43").Value = '0.870'
$ws.Range("E43").Value = '  -5.69%  '
$ws.Range("B44").Value = 'WhiteBITCoin'
$ws.Range("C44").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D44").Value = '24.13'
$ws.Range("E44").Value = '  -0.06%  '
$ws.Range("B45").Value = 'MantraDAO'
$ws.Range("C45").Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range("D45").Value = '3.74'
$ws.Range("E45").Value = '  +3.50%  '
$ws.Range("D46").Value = '0.0416'
$ws.Range("E46").Value = '  -1.11%  '
$ws.Range("D47").Value = '1.68'
$ws.Range("E47").Value = '  -2.73%  '
$ws.Range("D48").Value = '5.39'
$ws.Range("E48").Value = '  -3.76%  '
$ws.Range("D49").Value = '53.27'
$ws.Range("E49").Value = '  -0.65%  '
$ws.Range("D50").Value = '2.10'
$ws.Range("E50").Value = '  -4.25%  '
$ws.Range("D51").Value = '8.00'
$ws.Range("E51").Value = '  -2.27%  '

# Restore the default (unstyled) cell style so no stray formatting is
# introduced; only the text content of the cells should differ.
$dataRange.Style = "Normal"
